$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hazard-focused")

# The old "Hazard words" header (A1) is being replaced by two new headers:
# B1 = "Action/Descriptor", A1 = "Hazard Noun/Subject" (set B1 first so the
# shared-string table ends up in the same order as the target workbook).
$ws1.Range("B1").Value = "Action/Descriptor"
$ws1.Range("A1").Value = "Hazard Noun/Subject"

# Give the two new columns explicit widths, matching the authored column sizing.
$ws1.Columns.Item(1).ColumnWidth = 26.21875
$ws1.Columns.Item(2).ColumnWidth = 25.77734375

# Move the active selection from M1 to E1 on the Hazard-focused sheet.
$ws1.Range("E1").Select() | Out-Null
